$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = '39.412.58'
$ws.Range("E2").Value = '  +1.77%  '

# Row 3 - Ethereum
$ws.Range("D3").Value = '2.167.17'
$ws.Range("E3").Value = '  +3.62%  '

# Row 4 - TetherUSD
$ws.Range("E4").Value = '  +0.09%  '

# Row 5 - BNB
$ws.Range("D5").Value = '''229.56'
$ws.Range("E5").Value = '  +0.16%  '

# Row 6 - XRP
$ws.Range("E6").Value = '  +1.16%  '

# Row 7 - Solana
$ws.Range("D7").Value = '''64.98'
$ws.Range("E7").Value = '  +6.28%  '

# Row 8 - USDC
$ws.Range("E8").Value = '  +0.08%  '

# Row 9 - Cardano
$ws.Range("E9").Value = '  +3.60%  '

# Row 10 - Dogecoin
$ws.Range("D10").Value = '''0.0862'
$ws.Range("E10").Value = '  +2.03%  '

# Row 11 - TRON
$ws.Range("E11").Value = '  -0.15%  '

# Row 12 - Chainlink
$ws.Range("D12").Value = '''15.95'
$ws.Range("E12").Value = '  +4.84%  '

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = '2.491.44'
$ws.Range("E13").Value = '  +3.79%  '

# Row 14 - Avalanche
$ws.Range("D14").Value = '''22.43'
$ws.Range("E14").Value = '  +2.07%  '

# Row 15 - Polygon
$ws.Range("D15").Value = '''0.815'
$ws.Range("E15").Value = '  +0.30%  '

# Row 16 - Polkadot
$ws.Range("D16").Value = '''5.57'
$ws.Range("E16").Value = '  +2.20%  '

# Row 17 - WrappedEther
$ws.Range("D17").Value = '2.153.57'
$ws.Range("E17").Value = '  +2.85%  '

# Row 18 - WrappedBTC
$ws.Range("D18").Value = '39.441.28'
$ws.Range("E18").Value = '  +2.01%  '

# Row 19 - Uniswap->Litecoin
$ws.Range("B19").Value = 'Litecoin'
$ws.Range("C19").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D19").Value = '''72.29'
$ws.Range("E19").Value = '  +0.73%  '

# Row 20 - Litecoin->Uniswap
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = '''6.17'
$ws.Range("E20").Value = '  +1.29%  '

# Row 21 - ShibaInu
$ws.Range("D21").Value = '0.0₃0855'
$ws.Range("E21").Value = '  +1.51%  '

# Row 22 - BitcoinCash
$ws.Range("D22").Value = '''232.19'
$ws.Range("E22").Value = '  +2.04%  '

# Row 24 - Toncoin
$ws.Range("D24").Value = '''2.41'
$ws.Range("E24").Value = '  +0.98%  '

# Row 25 - PancakeSwap
$ws.Range("D25").Value = '''2.38'
$ws.Range("E25").Value = '  +2.39%  '

# Row 26 - Monero->Cosmos
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").Value = '''9.57'
$ws.Range("E26").Value = '  +0.39%  '

# Row 27 - Cosmos->Monero
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = '''172.42'
$ws.Range("E27").Value = '  +0.94%  '

# Row 28 - Kaspa
$ws.Range("E28").Value = '  -0.65%  '

# Row 29 - EthereumClassic
$ws.Range("D29").Value = '''20.03'
$ws.Range("E29").Value = '  +4.08%  '

# Row 30 - ImmutableX
$ws.Range("D30").Value = '''1.40'
$ws.Range("E30").Value = '  -2.17%  '

# Row 31 - WEMIXToken
$ws.Range("E31").Value = '  +11.92%  '

# Row 32 - Stellar
$ws.Range("E32").Value = '  +1.64%  '

# Row 33 - Filecoin
$ws.Range("D33").Value = '''4.65'
$ws.Range("E33").Value = '  +3.06%  '

# Row 34 - InternetComputer(DFINITY)
$ws.Range("D34").Value = '''4.79'
$ws.Range("E34").Value = '  +2.16%  '

# Row 35 - THORChain
$ws.Range("D35").Value = '''7.13'
$ws.Range("E35").Value = '  +8.90%  '

# Row 36 - Hedera
$ws.Range("D36").Value = '''0.0620'
$ws.Range("E36").Value = '  +1.66%  '

# Row 37 - LidoDAOToken
$ws.Range("D37").Value = '''2.44'
$ws.Range("E37").Value = '  +1.65%  '

# Row 38 - RenderToken
$ws.Range("D38").Value = '''3.59'
$ws.Range("E38").Value = '  +0.35%  '

# Row 39 - BinanceUSD
$ws.Range("E39").Value = '  +0.09%  '

# Row 40 - Aave
$ws.Range("D40").Value = '''104.51'
$ws.Range("E40").Value = '  +3.40%  '

# Row 41 - VeChain
$ws.Range("D41").Value = '''0.0230'
$ws.Range("E41").Value = '  +0.53%  '

# Row 42 - InjectiveProtocol
$ws.Range("D42").Value = '''17.85'
$ws.Range("E42").Value = '  -0.92%  '

# Row 43 - Maker
$ws.Range("D43").Value = '1.540.97'
$ws.Range("E43").Value = '  +0.48%  '

# Row 44 - TrustWalletToken
$ws.Range("E44").Value = '  +5.08%  '

# Row 45 - FTXToken
$ws.Range("E45").Value = '  +7.10%  '

# Row 46 - FraxShare->ARBITRUM
$ws.Range("B46").Value = 'ARBITRUM'
$ws.Range("C46").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D46").Value = '''1.11'
$ws.Range("E46").Value = '  +7.17%  '

# Row 47 - ARBITRUM->Cronos
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").Value = '''0.0927'
$ws.Range("E47").Value = '  +1.10%  '

# Row 48 - Cronos->HuobiToken
$ws.Range("B48").Value = 'HuobiToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D48").Value = '''2.82'
$ws.Range("E48").Value = '  +0.62%  '

# Row 49 - HuobiToken->FraxShare
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").Value = '''7.87'
$ws.Range("E49").Value = '  +2.01%  '

# Row 50 - RocketPoolETH
$ws.Range("D50").Value = '2.374.96'
$ws.Range("E50").Value = '  +3.86%  '

# Row 51 - MXToken
$ws.Range("E51").Value = '  +0.28%  '
